$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Row 3 tweak: B3 "student" -> "d"
# ------------------------------------------------------------------
$ws.Range("B3").Value = "d"

# ------------------------------------------------------------------
# 2. Build out row 4 as a new guest-checkout record (mirrors rows 2/3)
# ------------------------------------------------------------------

# Start by cloning row 3's formatting onto row 4 so every column picks
# up the same cell styles (borders, hyperlink look, number formats...).
$ws.Range("A3:AS3").Copy()
$ws.Range("A4:AS4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
# Row 3 has no content in AJ:AN, but the format paste still stamps those
# cells in; drop them so row 4 stays empty there too (nothing in the
# source row used those columns).
$ws.Range("AJ4:AN4").Clear()

$ws.Rows(4).RowHeight = 57.6

# Values are written in the same left-to-right / field order the
# original author used, so brand-new shared-string entries land at the
# same index the recorded workbook expects (search, profile, viewClass,
# email, firstname, dependent_name, dept_email, ...).
$ws.Range("A4").Value = "guestckt"
$ws.Range("B4").Value = "student"
$ws.Range("C4").Value = "Painting"
$ws.Range("D4").Value = "mgs"
$ws.Range("E4").Value = "gms-academy/paint-sess-1"
$ws.Range("H4").Value = "sri3@nkt12.com"
$ws.Range("F4").Value = "ranku"
$ws.Range("G4").Value = 94589184002
$ws.Range("I4").Value = "Test@1234"
$ws.Range("J4").Value = "Test@1234"
$ws.Range("K4").Value = "sakasuki"
$ws.Range("L4").Value = 21
$ws.Range("M4").Value = "male"
$ws.Range("N4").Value = "cousin"
$ws.Range("O4").Value = "rai17@nkt6.com"
$ws.Range("P4").Value = "Test@1234"
$ws.Range("Q4").Value = "Test@1234"
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = "Wednesday"
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 4
$ws.Range("X4").Formula = "=W4+AA4"
$ws.Range("Y4").Formula = "=(((V4+1)*X4)*R4)"
$ws.Range("AA4").Formula = "=W4*AB4"
$ws.Range("AB4").Value = 0.25
$ws.Range("AC4").Value = 0.01
$ws.Range("AD4").Formula = "=Y4*AC4"
$ws.Range("AE4").Formula = "=AC4*Y4"
$ws.Range("AF4").Formula = "=Y4+AE4"
$ws.Range("AI4").Value = "FLAT53"
$ws.Range("AO4").Value = "4242 4242 4242 4242"
$ws.Range("AP4").Value = "12/35"
$ws.Range("AQ4").Value = "4580"
$ws.Range("AR4").Value = 1
$ws.Range("AS4").Value = "Login Success > Class booked"

# Hyperlinks on row 4 (same columns as rows 2 & 3)
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:sri3@nkt12.com")
$ws.Hyperlinks.Add($ws.Range("I4"), "mailto:Test@1234")
$ws.Hyperlinks.Add($ws.Range("J4"), "mailto:Test@1234")
$ws.Hyperlinks.Add($ws.Range("O4"), "mailto:rai17@nkt6.com")
$ws.Hyperlinks.Add($ws.Range("P4"), "mailto:Test@1234")
$ws.Hyperlinks.Add($ws.Range("Q4"), "mailto:Test@1234")

# Re-stamp the formatting so the Hyperlinks.Add styling doesn't clobber
# the cell styles we copied from row 3 above.
$ws.Range("A3:AS3").Copy()
$ws.Range("A4:AS4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("AJ4:AN4").Clear()
$ws.Rows(4).RowHeight = 57.6

# ------------------------------------------------------------------
# 3. Selection moves to O4 (matches the saved cursor position)
# ------------------------------------------------------------------
$ws.Range("O4").Select()
